$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the latex_name column (column B) to wrap the LaTeX math segments
# with '$' delimiters, moving the formatting from code into the spreadsheet.
$ws.Range("B2").Value = '$\delta ^{18}$O'
$ws.Range("B3").Value = '$\delta$ D'
$ws.Range("B4").Value = 'd$_{xs}$'
$ws.Range("B5").Value = 'CO$_{2}$'
$ws.Range("B6").Value = 'CH$_{4}$'

# Update the active selection to match the saved view state.
$ws.Range("B6").Select()
